$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.262426
$ws.Range("H2").Value = 0.7872779999999999
$ws.Range("I2").Value = 0.2164351829578579
$ws.Range("J2").Value = 0.2164351829578578
$ws.Range("M2").Value = 0.05993299999999999
$ws.Range("Q2").Value = 0.015727977458
$ws.Range("R2").Value = 0.141551797122
$ws.Range("S2").Value = 0.2164351829578579
$ws.Range("T2").Value = 0.2164351829578578

# Row 3 (FAPs)
$ws.Range("I3").Value = 0.3553218343373718
$ws.Range("J3").Value = 0.3553218343373717
$ws.Range("M3").Value = 0.05993299999999999
$ws.Range("S3").Value = 0.3553218343373718
$ws.Range("T3").Value = 0.3553218343373717

# Row 4 (MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1998913333333333
$ws.Range("H4").Value = 0.599674
$ws.Range("I4").Value = 0.1648598740280695
$ws.Range("J4").Value = 0.1648598740280695
$ws.Range("M4").Value = 0.05993299999999999
$ws.Range("Q4").Value = 0.01198008728066667
$ws.Range("R4").Value = 0.107820785526
$ws.Range("S4").Value = 0.1648598740280695
$ws.Range("T4").Value = 0.1648598740280695

# Row 5 (Resolving-Mac)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.31935
$ws.Range("H5").Value = 0.95805
$ws.Range("I5").Value = 0.2633831086767009
$ws.Range("J5").Value = 0.2633831086767008
$ws.Range("M5").Value = 0.05993299999999999
$ws.Range("Q5").Value = 0.01913960355
$ws.Range("R5").Value = 0.17225643195
$ws.Range("S5").Value = 0.2633831086767009
$ws.Range("T5").Value = 0.2633831086767008
